# Fill the "Personel" sheet (2nd tab) with two names and the "Isler"
# sheet (3rd tab) with three category labels, then leave the "Isler"
# tab as the active/selected one - mirroring a user typing data into
# both sheets and finishing on the Isler tab.

$wb = $excel.ActiveWorkbook

# --- Personel sheet ---------------------------------------------------
$wsPersonel = $wb.Worksheets.Item(2)
$wsPersonel.Range("A1").Value = "Gökhan ELGÜL"
$wsPersonel.Range("A2").Value = "Göktan ELGÜL"
[void]$wsPersonel.Columns.Item(1).AutoFit()
[void]$wsPersonel.Range("A3").Select()

# --- Isler sheet --------------------------------------------------------
$wsIsler = $wb.Worksheets.Item(3)
$wsIsler.Range("A1").Value = "Cins D."
$wsIsler.Range("A2").Value = "İfraz"
$wsIsler.Range("A3").Value = "Tevhid"
[void]$wsIsler.Range("A4").Select()

# Isler ends up the active tab
[void]$wsIsler.Activate()
